# Processed Results - updated measures
# Adds a "Std" / "Relative std" statistic pair to each of the three
# CPU-frequency blocks (Low / Medium / High) on the worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each block: label row (Std / Relative std headers) then a formula row
# directly below it (STDEV + relative STDEV as a percentage of the mean).
$blocks = @(
    @{ LabelRow = 15; FormulaRow = 16; Range = "B2:B31";  MeanCell = "E4"  },
    @{ LabelRow = 47; FormulaRow = 48; Range = "B34:B63"; MeanCell = "E36" },
    @{ LabelRow = 79; FormulaRow = 80; Range = "B66:B95"; MeanCell = "E68" }
)

foreach ($block in $blocks) {
    $labelRow = $block.LabelRow
    $formulaRow = $block.FormulaRow
    $range = $block.Range
    $meanCell = $block.MeanCell

    $dLabel = $ws.Range("D$labelRow")
    $eLabel = $ws.Range("E$labelRow")
    $dLabel.Value = "Std"
    $eLabel.Value = "Relative std"
    $ws.Range("D${labelRow}:E${labelRow}").Font.Bold = $true

    $ws.Range("D$formulaRow").Formula = "=STDEV($range)"
    $ws.Range("E$formulaRow").Formula = "=(D$formulaRow/$meanCell)*100"
}

# Update the saved view state: scroll back to the top and move the
# active selection to E48 (matches the latest edit location).
$ws.Range("E48").Select()
